$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("SearchForEmployees")
$ws5 = $wb.Worksheets.Item("AddEmployee")

# Update SearchForEmployees data values (replace retired test-data strings
# with the new ones). Order matters so the rebuilt shared-string table ends
# up with the new strings appended in this sequence.
$ws1.Range("C6").Value = "Consultant"
$ws1.Range("C9").Value = "Greensboro"
$ws1.Range("C8").Value = "Bluesource"
$ws1.Range("C3").Value = "Kristi"
$ws1.Range("C5").Value = "Kevin"
$ws1.Range("C7").Value = "Perry Thomas"

# Move the selection/active-cell on AddEmployee first (so it no longer is
# the last-activated sheet), then finish on SearchForEmployees so it ends
# up as the workbook's active tab with C8 selected.
$ws5.Range("B8").Select() | Out-Null
$ws1.Range("C8").Select() | Out-Null
